# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2..21
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 2
    20 = 0
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
